$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 158
$ws.Range("I2").Value = 430
$ws.Range("J2").Value = 1876
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 554
$ws.Range("M2").Value = 28
$ws.Range("N2").Value = 312
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 25
$ws.Range("S2").Value = 192
$ws.Range("U2").Value = 18
$ws.Range("V2").Value = 2824
$ws.Range("X2").Value = 2834
$ws.Range("Z2").Value = 40
$ws.Range("AA2").Value = 26
